$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers I1 and J1, matching the style of the existing header H1
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Data for columns I (I0) and J (IF), rows 2-38
$data = @{
    2  = @(1, 4)
    3  = @(1, 2)
    4  = @(4, 4)
    5  = @(7, 9)
    6  = @(1, 4)
    7  = @(1, 6)
    8  = @(1, 5)
    9  = @(1, 6)
    10 = @(1, 6)
    11 = @(1, 7)
    12 = @(1, 6)
    13 = @(1, 5)
    14 = @(1, 7)
    15 = @(1, 6)
    16 = @(1, 4)
    17 = @(1, 8)
    18 = @(1, 5)
    19 = @(1, 8)
    20 = @(1, 6)
    21 = @(1, 3)
    22 = @(1, 9)
    23 = @(1, 8)
    24 = @(1, 6)
    25 = @(1, 8)
    26 = @(1, 7)
    27 = @(1, 5)
    28 = @(1, 4)
    29 = @(1, 4)
    30 = @(1, 7)
    31 = @(1, 8)
    32 = @(1, 5)
    33 = @(1, 7)
    34 = @(1, 6)
    35 = @(1, 5)
    36 = @(1, 5)
    37 = @(1, 3)
    38 = @(1, 1)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
